# TIC-8 Ticketplan FTP documentation (#196)
# Adds a new "policy.data.order_currency" column (AF) to the Ticketplan
# import template, populated with "euro" for both sample rows, and moves
# the sheet's viewport/selection over to the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + sample values in column AF (32nd column)
$ws.Range("AF1").Value = "policy.data.order_currency"
$ws.Range("AF2").Value = "euro"
$ws.Range("AF3").Value = "euro"

# Match the column width used for the new column in the diff
$ws.Range("AF1").ColumnWidth = 21.6640625

# Update viewport/selection to show the newly added column
$ws.Range("AA1").Select()
$ws.Range("AD15").Select()
